$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2602.3442
$ws.Range("I137").Value = 2475.3022
$ws.Range("J137").Value = 2905.8333
$ws.Range("K137").Value = 7425.9066
$ws.Range("L137").Value = 8717.499899999999
$ws.Range("M137").Value = -4875.9066
$ws.Range("N137").Value = -13817.4999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3707387.5
$ws.Range("I32").Value = 3165.7532
$ws.Range("K32").Value = 3165.7532
$ws.Range("M32").Value = -2878.7532
$ws.Range("H74").Value = 2383.3333
$ws.Range("I74").Value = 1583.3334
$ws.Range("J74").Value = 3183.3333
$ws.Range("K74").Value = 1583.3334
$ws.Range("L74").Value = 3183.3333
$ws.Range("M74").Value = -709.3334
$ws.Range("N74").Value = -4931.3333
$ws.Range("H77").Value = 2383.3333
$ws.Range("I77").Value = 1583.3334
$ws.Range("J77").Value = 3183.3333
$ws.Range("K77").Value = 7916.666999999999
$ws.Range("L77").Value = 15916.6665
$ws.Range("M77").Value = -3548.666999999999
$ws.Range("N77").Value = -24652.6665
$ws.Range("H132").Value = 701677.5600000001
$ws.Range("I132").Value = 1259.5231
$ws.Range("J132").Value = 3097844.5
$ws.Range("K132").Value = 3778.5693
$ws.Range("L132").Value = 9293533.5
$ws.Range("M132").Value = -1248.5693
$ws.Range("N132").Value = -9298593.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1595.6
$ws.Range("I16").Value = 448.4
$ws.Range("J16").Value = 2169.2
$ws.Range("K16").Value = 448.4
$ws.Range("L16").Value = 2169.2
$ws.Range("M16").Value = -161.4
$ws.Range("N16").Value = -2743.2
$ws.Range("H31").Value = 1195.6102
$ws.Range("I31").Value = 905.55554
$ws.Range("J31").Value = 1322.9512
$ws.Range("K31").Value = 905.55554
$ws.Range("L31").Value = 1322.9512
$ws.Range("M31").Value = -610.55554
$ws.Range("N31").Value = -1912.9512
$ws.Range("H34").Value = 1195.6102
$ws.Range("I34").Value = 905.55554
$ws.Range("J34").Value = 1322.9512
$ws.Range("K34").Value = 905.55554
$ws.Range("L34").Value = 1322.9512
$ws.Range("M34").Value = -703.55554
$ws.Range("N34").Value = -1726.9512
$ws.Range("H58").Value = 32258988
$ws.Range("I58").Value = 41667444
$ws.Range("J58").Value = 1428.4286
$ws.Range("K58").Value = 41667444
$ws.Range("L58").Value = 1428.4286
$ws.Range("M58").Value = -41667241
$ws.Range("N58").Value = -1834.4286
$ws.Range("H62").Value = 4631.25
$ws.Range("I62").Value = 2866.3333
$ws.Range("J62").Value = 5219.5557
$ws.Range("K62").Value = 2866.3333
$ws.Range("L62").Value = 5219.5557
$ws.Range("M62").Value = -2242.3333
$ws.Range("N62").Value = -6467.5557
$ws.Range("H65").Value = 4631.25
$ws.Range("I65").Value = 2866.3333
$ws.Range("J65").Value = 5219.5557
$ws.Range("K65").Value = 14331.6665
$ws.Range("L65").Value = 26097.7785
$ws.Range("M65").Value = -11211.6665
$ws.Range("N65").Value = -32337.7785
$ws.Range("H113").Value = 1595.6
$ws.Range("I113").Value = 448.4
$ws.Range("J113").Value = 2169.2
$ws.Range("K113").Value = 448.4
$ws.Range("L113").Value = 2169.2
$ws.Range("M113").Value = 1721.6
$ws.Range("N113").Value = -6509.2
$ws.Range("H136").Value = 32258988
$ws.Range("I136").Value = 41667444
$ws.Range("J136").Value = 1428.4286
$ws.Range("K136").Value = 125002332
$ws.Range("L136").Value = 4285.2858
$ws.Range("M136").Value = -124999782
$ws.Range("N136").Value = -9385.2858

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H51").Value = 1014.7143
$ws.Range("I51").Value = 368
$ws.Range("J51").Value = 1499.75
$ws.Range("K51").Value = 1104
$ws.Range("L51").Value = 4499.25
$ws.Range("M51").Value = -644
$ws.Range("N51").Value = -5419.25
$ws.Range("H121").Value = 44444890
$ws.Range("I121").Value = 319.75
$ws.Range("J121").Value = 80000540
$ws.Range("K121").Value = 959.25
$ws.Range("L121").Value = 240001620
$ws.Range("M121").Value = 350.75
$ws.Range("N121").Value = -240004240

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H69").Value = 19800
$ws.Range("J69").Value = 19800
$ws.Range("L69").Value = 19800
$ws.Range("N69").Value = -21298
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H72").Value = 19800
$ws.Range("J72").Value = 19800
$ws.Range("L72").Value = 59400
$ws.Range("N72").Value = -66888
$ws.Range("H74").Value = 58866.668
$ws.Range("J74").Value = 58866.668
$ws.Range("L74").Value = 58866.668
$ws.Range("N74").Value = -60738.668
$ws.Range("H77").Value = 58866.668
$ws.Range("J77").Value = 58866.668
$ws.Range("L77").Value = 176600.004
$ws.Range("N77").Value = -185960.004
$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("M80").ClearContents()
$ws.Range("N80").ClearContents()
$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("M83").ClearContents()
$ws.Range("N83").ClearContents()
$ws.Range("H113").Value = 1825
$ws.Range("J113").Value = 2450
$ws.Range("L113").Value = 2450
$ws.Range("N113").Value = -6790

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 14660.9
$ws.Range("J132").Value = 14660.9
$ws.Range("L132").Value = 43982.7
$ws.Range("N132").Value = -49042.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 628.3333
$ws.Range("I100").Value = 539
$ws.Range("J100").Value = 717.6667
$ws.Range("K100").Value = 1078
$ws.Range("L100").Value = 1435.3334
$ws.Range("M100").Value = -537
$ws.Range("N100").Value = -2517.3334
$ws.Range("H132").Value = 12376.637
$ws.Range("J132").Value = 12376.637
$ws.Range("L132").Value = 37129.911
$ws.Range("N132").Value = -42189.911
